# Actualización automática 2025-07-21 14:40:09
# Updates the July sale amount for client "JARAMILLO CARVAJAL NICOLAS ESTEBAN"
# (advisor HIDALGO HIDALGO PEDRO GUSTAVO) in the "NO RESURTIBLES" group,
# and propagates the change through the dependent summary sheets.

$wb = $excel.ActiveWorkbook

$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Sheet "VENTAS POR GRUPO" ---
# P10: NO RESURTIBLES sale for JARAMILLO CARVAJAL NICOLAS ESTEBAN
$wsGrupo.Range("P10").Value = 110.13
# P22: recalculated "non-zero count out of 20" label for column P
$wsGrupo.Range("P22").Value = "1 de 20"

# --- Sheet "VENTA MENSUAL" ---
# F10: julio sale for JARAMILLO CARVAJAL NICOLAS ESTEBAN
$wsMensual.Range("F10").Value = 6415.05
# F22: julio total
$wsMensual.Range("F22").Value = 36450.05

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
# Row 10: NO RESURTIBLES group totals
$wsCumplimiento.Range("D10").Value = 110.13
$wsCumplimiento.Range("E10").Value = 540.12
$wsCumplimiento.Range("F10").Value = 0.1693656286043829

# Row 19: TOTAL row
$wsCumplimiento.Range("D19").Value = 36450.05
$wsCumplimiento.Range("E19").Value = 28927.94762291769
$wsCumplimiento.Range("F19").Value = 0.5575277819035369
